$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellA = $ws.Cells.Item(31, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "09/17/2025"

$ws.Cells.Item(31, 2).Value = 15900.13
